# Update "想去人数" (number of people interested) counts in the
# "展览" and "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F (想去人数)
$updates = @{
    2  = 624
    3  = 2191
    4  = 80
    5  = 12961
    6  = 69
    9  = 475
    10 = 1167
    11 = 970
    12 = 13732
    13 = 14238
    22 = 1083
    26 = 5319
    27 = 11
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
